# Populate the previously-empty "SeptemberRaw" sheet with the raw
# per-library sharing statistics for September 2023. Everything that
# depends on it ("September" summary sheet and "Yearly total") pulls
# these numbers in through existing formulas and recalculates automatically.

$wb = $excel.ActiveWorkbook
$originalActiveSheet = $wb.ActiveSheet.Name

$ws = $wb.Worksheets.Item("SeptemberRaw")

# Header row
$ws.Cells.Item(1, 1).Value = 'Library'
$ws.Cells.Item(1, 2).Value = 'Items owned by this library checked out at this library this month'
$ws.Cells.Item(1, 3).Value = 'Items owned by other libraries checked out at this library this month'
$ws.Cells.Item(1, 4).Value = 'Total circulation this month'

# Row 2: 'Atchison Public Library'
$ws.Cells.Item(2, 1).Value = 'Atchison Public Library'
$ws.Cells.Item(2, 2).Value = 4102
$ws.Cells.Item(2, 3).Value = 1586
$ws.Cells.Item(2, 4).Value = 5688

# Row 3: 'Baldwin City Public Library'
$ws.Cells.Item(3, 1).Value = 'Baldwin City Public Library'
$ws.Cells.Item(3, 2).Value = 2780
$ws.Cells.Item(3, 3).Value = 820
$ws.Cells.Item(3, 4).Value = 3600

# Row 4: 'Basehor Community Library'
$ws.Cells.Item(4, 1).Value = 'Basehor Community Library'
$ws.Cells.Item(4, 2).Value = 8280
$ws.Cells.Item(4, 3).Value = 1263
$ws.Cells.Item(4, 4).Value = 9543

# Row 5: 'Bern Community Library'
$ws.Cells.Item(5, 1).Value = 'Bern Community Library'
$ws.Cells.Item(5, 2).Value = 147
$ws.Cells.Item(5, 3).Value = 124
$ws.Cells.Item(5, 4).Value = 271

# Row 6: 'Bonner Springs City Library'
$ws.Cells.Item(6, 1).Value = 'Bonner Springs City Library'
$ws.Cells.Item(6, 2).Value = 5392
$ws.Cells.Item(6, 3).Value = 1315
$ws.Cells.Item(6, 4).Value = 6707

# Row 7: 'Burlingame Community Library'
$ws.Cells.Item(7, 1).Value = 'Burlingame Community Library'
$ws.Cells.Item(7, 2).Value = 398
$ws.Cells.Item(7, 3).Value = 172
$ws.Cells.Item(7, 4).Value = 570

# Row 8: 'Carbondale City Library'
$ws.Cells.Item(8, 1).Value = 'Carbondale City Library'
$ws.Cells.Item(8, 2).Value = 664
$ws.Cells.Item(8, 3).Value = 163
$ws.Cells.Item(8, 4).Value = 827

# Row 9: 'Centralia Community Library'
$ws.Cells.Item(9, 1).Value = 'Centralia Community Library'
$ws.Cells.Item(9, 2).Value = 273
$ws.Cells.Item(9, 3).Value = 37
$ws.Cells.Item(9, 4).Value = 310

# Row 10: 'Corning City Library'
$ws.Cells.Item(10, 1).Value = 'Corning City Library'
$ws.Cells.Item(10, 2).Value = 6
$ws.Cells.Item(10, 4).Value = 6

# Row 11: 'Digital Content'
$ws.Cells.Item(11, 1).Value = 'Digital Content'

# Row 12: 'Doniphan County Library - Elwood'
$ws.Cells.Item(12, 1).Value = 'Doniphan County Library - Elwood'
$ws.Cells.Item(12, 2).Value = 127
$ws.Cells.Item(12, 3).Value = 11
$ws.Cells.Item(12, 4).Value = 138

# Row 13: 'Doniphan County Library - Highland'
$ws.Cells.Item(13, 1).Value = 'Doniphan County Library - Highland'
$ws.Cells.Item(13, 2).Value = 256
$ws.Cells.Item(13, 3).Value = 148
$ws.Cells.Item(13, 4).Value = 404

# Row 14: 'Doniphan County Library - Troy'
$ws.Cells.Item(14, 1).Value = 'Doniphan County Library - Troy'
$ws.Cells.Item(14, 2).Value = 462
$ws.Cells.Item(14, 3).Value = 108
$ws.Cells.Item(14, 4).Value = 570

# Row 15: 'Doniphan County Library - Wathena'
$ws.Cells.Item(15, 1).Value = 'Doniphan County Library - Wathena'
$ws.Cells.Item(15, 2).Value = 446
$ws.Cells.Item(15, 3).Value = 127
$ws.Cells.Item(15, 4).Value = 573

# Row 16: 'Effingham Community Library'
$ws.Cells.Item(16, 1).Value = 'Effingham Community Library'
$ws.Cells.Item(16, 2).Value = 238
$ws.Cells.Item(16, 3).Value = 35
$ws.Cells.Item(16, 4).Value = 273

# Row 17: 'Eudora Community Library'
$ws.Cells.Item(17, 1).Value = 'Eudora Community Library'
$ws.Cells.Item(17, 2).Value = 1627
$ws.Cells.Item(17, 3).Value = 741
$ws.Cells.Item(17, 4).Value = 2368

# Row 18: 'Everest, Barnes Reading Room'
$ws.Cells.Item(18, 1).Value = 'Everest, Barnes Reading Room'
$ws.Cells.Item(18, 2).Value = 207
$ws.Cells.Item(18, 3).Value = 88
$ws.Cells.Item(18, 4).Value = 295

# Row 19: 'Hiawatha, Morrill Public Library'
$ws.Cells.Item(19, 1).Value = 'Hiawatha, Morrill Public Library'
$ws.Cells.Item(19, 2).Value = 1630
$ws.Cells.Item(19, 3).Value = 629
$ws.Cells.Item(19, 4).Value = 2259

# Row 20: 'Highland Community College'
$ws.Cells.Item(20, 1).Value = 'Highland Community College'
$ws.Cells.Item(20, 2).Value = 133
$ws.Cells.Item(20, 3).Value = 4
$ws.Cells.Item(20, 4).Value = 137

# Row 21: 'Holton, Beck-Bookman Library'
$ws.Cells.Item(21, 1).Value = 'Holton, Beck-Bookman Library'
$ws.Cells.Item(21, 2).Value = 1835
$ws.Cells.Item(21, 3).Value = 436
$ws.Cells.Item(21, 4).Value = 2271

# Row 22: 'Horton Public Library'
$ws.Cells.Item(22, 1).Value = 'Horton Public Library'
$ws.Cells.Item(22, 2).Value = 101
$ws.Cells.Item(22, 3).Value = 38
$ws.Cells.Item(22, 4).Value = 139

# Row 23: 'Lansing Community Library'
$ws.Cells.Item(23, 1).Value = 'Lansing Community Library'
$ws.Cells.Item(23, 2).Value = 2207
$ws.Cells.Item(23, 3).Value = 724
$ws.Cells.Item(23, 4).Value = 2931

# Row 24: 'Leavenworth Public Library'
$ws.Cells.Item(24, 1).Value = 'Leavenworth Public Library'
$ws.Cells.Item(24, 2).Value = 8772
$ws.Cells.Item(24, 3).Value = 1997
$ws.Cells.Item(24, 4).Value = 10769

# Row 25: 'Linwood Community Library'
$ws.Cells.Item(25, 1).Value = 'Linwood Community Library'
$ws.Cells.Item(25, 2).Value = 626
$ws.Cells.Item(25, 3).Value = 176
$ws.Cells.Item(25, 4).Value = 802

# Row 26: 'Louisburg Library'
$ws.Cells.Item(26, 1).Value = 'Louisburg Library'

# Row 27: 'Lyndon Carnegie Library'
$ws.Cells.Item(27, 1).Value = 'Lyndon Carnegie Library'
$ws.Cells.Item(27, 2).Value = 455
$ws.Cells.Item(27, 3).Value = 176
$ws.Cells.Item(27, 4).Value = 631

# Row 28: 'McLouth Public Library'
$ws.Cells.Item(28, 1).Value = 'McLouth Public Library'
$ws.Cells.Item(28, 2).Value = 84
$ws.Cells.Item(28, 3).Value = 51
$ws.Cells.Item(28, 4).Value = 135

# Row 29: 'Meriden-Ozawkie Public Library'
$ws.Cells.Item(29, 1).Value = 'Meriden-Ozawkie Public Library'
$ws.Cells.Item(29, 2).Value = 1258
$ws.Cells.Item(29, 3).Value = 563
$ws.Cells.Item(29, 4).Value = 1821

# Row 30: 'Northeast Kansas Library System'
$ws.Cells.Item(30, 1).Value = 'Northeast Kansas Library System'
$ws.Cells.Item(30, 2).Value = 4
$ws.Cells.Item(30, 3).Value = 39
$ws.Cells.Item(30, 4).Value = 43

# Row 31: 'Nortonville Public Library'
$ws.Cells.Item(31, 1).Value = 'Nortonville Public Library'
$ws.Cells.Item(31, 2).Value = 303
$ws.Cells.Item(31, 3).Value = 76
$ws.Cells.Item(31, 4).Value = 379

# Row 32: 'Osage City Library'
$ws.Cells.Item(32, 1).Value = 'Osage City Library'
$ws.Cells.Item(32, 2).Value = 1472
$ws.Cells.Item(32, 3).Value = 457
$ws.Cells.Item(32, 4).Value = 1929

# Row 33: 'Osawatomie Public Library'
$ws.Cells.Item(33, 1).Value = 'Osawatomie Public Library'
$ws.Cells.Item(33, 2).Value = 892
$ws.Cells.Item(33, 3).Value = 379
$ws.Cells.Item(33, 4).Value = 1271

# Row 34: 'Oskaloosa Public Library'
$ws.Cells.Item(34, 1).Value = 'Oskaloosa Public Library'
$ws.Cells.Item(34, 2).Value = 432
$ws.Cells.Item(34, 3).Value = 177
$ws.Cells.Item(34, 4).Value = 609

# Row 35: 'Ottawa Library'
$ws.Cells.Item(35, 1).Value = 'Ottawa Library'
$ws.Cells.Item(35, 2).Value = 6227
$ws.Cells.Item(35, 3).Value = 1060
$ws.Cells.Item(35, 4).Value = 7287

# Row 36: 'Overbrook Public Library'
$ws.Cells.Item(36, 1).Value = 'Overbrook Public Library'
$ws.Cells.Item(36, 2).Value = 794
$ws.Cells.Item(36, 3).Value = 179
$ws.Cells.Item(36, 4).Value = 973

# Row 37: 'Paola Free Library'
$ws.Cells.Item(37, 1).Value = 'Paola Free Library'
$ws.Cells.Item(37, 2).Value = 3043
$ws.Cells.Item(37, 3).Value = 549
$ws.Cells.Item(37, 4).Value = 3592

# Row 38: 'Perry-Lecompton Community Library'
$ws.Cells.Item(38, 1).Value = 'Perry-Lecompton Community Library'
$ws.Cells.Item(38, 2).Value = 79
$ws.Cells.Item(38, 3).Value = 20
$ws.Cells.Item(38, 4).Value = 99

# Row 39: 'Pomona Community Library'
$ws.Cells.Item(39, 1).Value = 'Pomona Community Library'
$ws.Cells.Item(39, 2).Value = 133
$ws.Cells.Item(39, 3).Value = 67
$ws.Cells.Item(39, 4).Value = 200

# Row 40: 'Prairie Hills Schools - Axtell Public School'
$ws.Cells.Item(40, 1).Value = 'Prairie Hills Schools - Axtell Public School'
$ws.Cells.Item(40, 2).Value = 643
$ws.Cells.Item(40, 3).Value = 49
$ws.Cells.Item(40, 4).Value = 692

# Row 41: 'Prairie Hills Schools - Sabetha Elementary School'
$ws.Cells.Item(41, 1).Value = 'Prairie Hills Schools - Sabetha Elementary School'
$ws.Cells.Item(41, 2).Value = 1907
$ws.Cells.Item(41, 3).Value = 91
$ws.Cells.Item(41, 4).Value = 1998

# Row 42: 'Prairie Hills Schools - Sabetha High School'
$ws.Cells.Item(42, 1).Value = 'Prairie Hills Schools - Sabetha High School'
$ws.Cells.Item(42, 2).Value = 25
$ws.Cells.Item(42, 3).Value = 5
$ws.Cells.Item(42, 4).Value = 30

# Row 43: 'Prairie Hills Schools - Sabetha Middle School'
$ws.Cells.Item(43, 1).Value = 'Prairie Hills Schools - Sabetha Middle School'
$ws.Cells.Item(43, 2).Value = 164
$ws.Cells.Item(43, 3).Value = 19
$ws.Cells.Item(43, 4).Value = 183

# Row 44: 'Prairie Hills Schools - Wetmore Academic Center (Permanently closed)'
$ws.Cells.Item(44, 1).Value = 'Prairie Hills Schools - Wetmore Academic Center (Permanently closed)'

# Row 45: 'Richmond Public Library'
$ws.Cells.Item(45, 1).Value = 'Richmond Public Library'
$ws.Cells.Item(45, 2).Value = 279
$ws.Cells.Item(45, 3).Value = 69
$ws.Cells.Item(45, 4).Value = 348

# Row 46: 'Rossville Community Library'
$ws.Cells.Item(46, 1).Value = 'Rossville Community Library'
$ws.Cells.Item(46, 2).Value = 1154
$ws.Cells.Item(46, 3).Value = 469
$ws.Cells.Item(46, 4).Value = 1623

# Row 47: 'Sabetha, Mary Cotton Library'
$ws.Cells.Item(47, 1).Value = 'Sabetha, Mary Cotton Library'
$ws.Cells.Item(47, 2).Value = 2963
$ws.Cells.Item(47, 3).Value = 892
$ws.Cells.Item(47, 4).Value = 3855

# Row 48: 'Seneca Free Library'
$ws.Cells.Item(48, 1).Value = 'Seneca Free Library'
$ws.Cells.Item(48, 2).Value = 1392
$ws.Cells.Item(48, 3).Value = 264
$ws.Cells.Item(48, 4).Value = 1656

# Row 49: 'Silver Lake Library'
$ws.Cells.Item(49, 1).Value = 'Silver Lake Library'
$ws.Cells.Item(49, 2).Value = 1032
$ws.Cells.Item(49, 3).Value = 422
$ws.Cells.Item(49, 4).Value = 1454

# Row 50: 'Tonganoxie Public Library'
$ws.Cells.Item(50, 1).Value = 'Tonganoxie Public Library'
$ws.Cells.Item(50, 2).Value = 3170
$ws.Cells.Item(50, 3).Value = 792
$ws.Cells.Item(50, 4).Value = 3962

# Row 51: 'Valley Falls, Delaware Township Library'
$ws.Cells.Item(51, 1).Value = 'Valley Falls, Delaware Township Library'
$ws.Cells.Item(51, 2).Value = 420
$ws.Cells.Item(51, 3).Value = 166
$ws.Cells.Item(51, 4).Value = 586

# Row 52: 'Wellsville City Library'
$ws.Cells.Item(52, 1).Value = 'Wellsville City Library'
$ws.Cells.Item(52, 2).Value = 1128
$ws.Cells.Item(52, 3).Value = 436
$ws.Cells.Item(52, 4).Value = 1564

# Row 53: 'Wetmore Public Library'
$ws.Cells.Item(53, 1).Value = 'Wetmore Public Library'
$ws.Cells.Item(53, 2).Value = 163
$ws.Cells.Item(53, 3).Value = 123
$ws.Cells.Item(53, 4).Value = 286

# Row 54: 'Williamsburg Community Library'
$ws.Cells.Item(54, 1).Value = 'Williamsburg Community Library'
$ws.Cells.Item(54, 2).Value = 362
$ws.Cells.Item(54, 3).Value = 19
$ws.Cells.Item(54, 4).Value = 381

# Row 55: 'Winchester Public Library'
$ws.Cells.Item(55, 1).Value = 'Winchester Public Library'
$ws.Cells.Item(55, 2).Value = 381
$ws.Cells.Item(55, 3).Value = 321
$ws.Cells.Item(55, 4).Value = 702

# The source sheet was saved without frozen panes/split selections; drop
# the inherited frozen-pane view state left over from the template sheet,
# then restore the workbook's originally active sheet.
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$wb.Worksheets.Item($originalActiveSheet).Activate()

Write-Host "SeptemberRaw populated; dependent sheets recalculated."
